# Weekly update: insert a new price record as row 20 (Camote, Vega Modelo
# de Temuco), pushing the existing rows 20-26 down to 21-27.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 20; rows 20-26 shift down to 21-27.
$ws.Rows.Item(20).Insert()

# Populate the newly inserted row 20 with the new record.
$ws.Cells.Item(20, 1).Value  = 10
$ws.Cells.Item(20, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(20, 3).Value  = "La Araucanía"
$ws.Cells.Item(20, 4).Value  = 44466
$ws.Cells.Item(20, 5).Value  = 9
$ws.Cells.Item(20, 6).Value  = 100114002
$ws.Cells.Item(20, 7).Value  = "Camote"
$ws.Cells.Item(20, 8).Value  = "Sin especificar"
$ws.Cells.Item(20, 9).Value  = "Primera"
$ws.Cells.Item(20, 10).Value = 20
$ws.Cells.Item(20, 11).Value = 25000
$ws.Cells.Item(20, 12).Value = 25000
$ws.Cells.Item(20, 13).Value = 25000
$ws.Cells.Item(20, 14).Value = "`$/malla 20 kilos"
$ws.Cells.Item(20, 15).Value = "Perú"
$ws.Cells.Item(20, 16).Value = 1667
$ws.Cells.Item(20, 17).Value = 15
$ws.Cells.Item(20, 18).Value = "Hortaliza"
